$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "558.04") are not auto-converted to floating point numbers,
# matching the original inline-string / text storage.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.121.89"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "2.418.06"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "558.04"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").Value = "143.19"
$ws.Range("E6").Value = "  +3.04%  "
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "2.414.29"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "26.19"
$ws.Range("E14").Value = "  +4.58%  "
$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").Value = "  +5.43%  "
$ws.Range("D16").Value = "2.848.49"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").Value = "61.890.87"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "2.415.54"
$ws.Range("D19").Value = "11.18"
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("D20").Value = "4.18"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "323.45"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "6.76"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "65.42"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "1.71"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "8.92"
$ws.Range("E26").Value = "  +6.39%  "
$ws.Range("D27").Value = "595.03"
$ws.Range("E27").Value = "  +17.38%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "2.520.95"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").Value = "0.0₃0937"
$ws.Range("E30").Value = "  +5.78%  "
$ws.Range("D31").Value = "8.27"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Value = "1.44"
$ws.Range("E32").Value = "  +4.91%  "
$ws.Range("D33").Value = "0.148"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").Value = "1.87"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("D36").Value = "5.68"
$ws.Range("E36").Value = "  +4.97%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "4.77"
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("D40").Value = "151.64"
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("D41").Value = "18.69"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "1.81"
$ws.Range("E42").Value = "  -4.99%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "2.35"
$ws.Range("E44").Value = "  +11.98%  "
$ws.Range("D45").Value = "150.68"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "3.65"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("D47").Value = "0.0539"
$ws.Range("E47").Value = "  +3.38%  "
$ws.Range("D48").Value = "20.22"
$ws.Range("E48").Value = "  +4.87%  "
$ws.Range("D49").Value = "0.591"
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("D50").Value = "0.0921"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").Value = "0.0229"
$ws.Range("E51").Value = "  +2.12%  "

# Restore original (default) cell formatting so no stray number-format styles
# are introduced into the saved workbook.
$ws.Range("D2:D51").ClearFormats()

